$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 75000
$ws.Cells.Item(15, 9).Value = 75000.5
$ws.Cells.Item(15, 10).Value = 0.1
$ws.Cells.Item(15, 11).Value = 985.2222
$ws.Cells.Item(15, 12).Value = 2955.6666
$ws.Cells.Item(15, 13).Value = -2786.6666
